$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.683228666666667
$ws.Range("H2").Value = 11.049686
$ws.Range("I2").Value = 0.1499797304438824
$ws.Range("J2").Value = 0.1499797304438824
$ws.Range("M2").Value = 1.050658666666667
$ws.Range("N2").Value = 3.151976
$ws.Range("O2").Value = 0.1514279810580986
$ws.Range("P2").Value = 0.1514279810580986
$ws.Range("Q2").Value = 3.869816119948446
$ws.Range("R2").Value = 34.82834507953601
$ws.Range("S2").Value = 0.02271112778075496
$ws.Range("T2").Value = 0.02271112778075496
$ws.Range("G3").Value = 3.683228666666667
$ws.Range("H3").Value = 11.049686
$ws.Range("I3").Value = 0.1499797304438824
$ws.Range("J3").Value = 0.1499797304438824
$ws.Range("M3").Value = 0.8776213333333334
$ws.Range("N3").Value = 2.632864
$ws.Range("O3").Value = 0.126488678822602
$ws.Range("P3").Value = 0.126488678822602
$ws.Range("Q3").Value = 3.232480053411556
$ws.Range("R3").Value = 29.092320480704
$ws.Range("S3").Value = 0.01897073795401666
$ws.Range("T3").Value = 0.01897073795401666
$ws.Range("G4").Value = 3.683228666666667
$ws.Range("H4").Value = 11.049686
$ws.Range("I4").Value = 0.1499797304438824
$ws.Range("J4").Value = 0.1499797304438824
$ws.Range("M4").Value = 5.010059000000001
$ws.Range("N4").Value = 15.030177
$ws.Range("O4").Value = 0.7220833401192995
$ws.Range("P4").Value = 0.7220833401192994
$ws.Range("Q4").Value = 18.45319293049134
$ws.Range("R4").Value = 166.078736374422
$ws.Range("S4").Value = 0.1082978647091108
$ws.Range("T4").Value = 0.1082978647091108
$ws.Range("I5").Value = 0.03800719241788433
$ws.Range("J5").Value = 0.03800719241788433
$ws.Range("M5").Value = 1.050658666666667
$ws.Range("N5").Value = 3.151976
$ws.Range("O5").Value = 0.1514279810580986
$ws.Range("P5").Value = 0.1514279810580986
$ws.Range("Q5").Value = 0.9806714911235557
$ws.Range("R5").Value = 8.826043420112001
$ws.Range("S5").Value = 0.005755352413526898
$ws.Range("T5").Value = 0.005755352413526897
$ws.Range("I6").Value = 0.03800719241788433
$ws.Range("J6").Value = 0.03800719241788433
$ws.Range("M6").Value = 0.8776213333333334
$ws.Range("N6").Value = 2.632864
$ws.Range("O6").Value = 0.126488678822602
$ws.Range("P6").Value = 0.126488678822602
$ws.Range("Q6").Value = 0.8191606359964444
$ws.Range("R6").Value = 7.372445723968
$ws.Range("S6").Value = 0.004807479554694604
$ws.Range("T6").Value = 0.004807479554694604
$ws.Range("I7").Value = 0.03800719241788433
$ws.Range("J7").Value = 0.03800719241788433
$ws.Range("M7").Value = 5.010059000000001
$ws.Range("N7").Value = 15.030177
$ws.Range("O7").Value = 0.7220833401192995
$ws.Range("P7").Value = 0.7220833401192994
$ws.Range("Q7").Value = 4.676325609852667
$ws.Range("R7").Value = 42.086930488674
$ws.Range("S7").Value = 0.02744436044966284
$ws.Range("T7").Value = 0.02744436044966283
$ws.Range("G8").Value = 4.395451333333333
$ws.Range("H8").Value = 13.186354
$ws.Range("I8").Value = 0.1789811781490995
$ws.Range("J8").Value = 0.1789811781490995
$ws.Range("M8").Value = 1.050658666666667
$ws.Range("N8").Value = 3.151976
$ws.Range("O8").Value = 0.1514279810580986
$ws.Range("P8").Value = 0.1514279810580986
$ws.Range("Q8").Value = 4.618119037278222
$ws.Range("R8").Value = 41.56307133550401
$ws.Range("S8").Value = 0.02710275845451801
$ws.Range("T8").Value = 0.027102758454518
$ws.Range("G9").Value = 4.395451333333333
$ws.Range("H9").Value = 13.186354
$ws.Range("I9").Value = 0.1789811781490995
$ws.Range("J9").Value = 0.1789811781490995
$ws.Range("M9").Value = 0.8776213333333334
$ws.Range("N9").Value = 2.632864
$ws.Range("O9").Value = 0.126488678822602
$ws.Range("P9").Value = 0.126488678822602
$ws.Range("Q9").Value = 3.857541859761778
$ws.Range("R9").Value = 34.717876737856
$ws.Range("S9").Value = 0.02263909275819235
$ws.Range("T9").Value = 0.02263909275819235
$ws.Range("G10").Value = 4.395451333333333
$ws.Range("H10").Value = 13.186354
$ws.Range("I10").Value = 0.1789811781490995
$ws.Range("J10").Value = 0.1789811781490995
$ws.Range("M10").Value = 5.010059000000001
$ws.Range("N10").Value = 15.030177
$ws.Range("O10").Value = 0.7220833401192995
$ws.Range("P10").Value = 0.7220833401192994
$ws.Range("Q10").Value = 22.02147051162867
$ws.Range("R10").Value = 198.193234604658
$ws.Range("S10").Value = 0.1292393269363891
$ws.Range("T10").Value = 0.1292393269363891
$ws.Range("G11").Value = 2.406480666666666
$ws.Range("H11").Value = 7.219442
$ws.Range("I11").Value = 0.09799101667823351
$ws.Range("J11").Value = 0.09799101667823353
$ws.Range("M11").Value = 1.050658666666667
$ws.Range("N11").Value = 3.151976
$ws.Range("O11").Value = 0.1514279810580986
$ws.Range("P11").Value = 0.1514279810580986
$ws.Range("Q11").Value = 2.528389768599111
$ws.Range("R11").Value = 22.755507917392
$ws.Range("S11").Value = 0.01483858181741537
$ws.Range("T11").Value = 0.01483858181741537
$ws.Range("G12").Value = 2.406480666666666
$ws.Range("H12").Value = 7.219442
$ws.Range("I12").Value = 0.09799101667823351
$ws.Range("J12").Value = 0.09799101667823353
$ws.Range("M12").Value = 0.8776213333333334
$ws.Range("N12").Value = 2.632864
$ws.Range("O12").Value = 0.126488678822602
$ws.Range("P12").Value = 0.126488678822602
$ws.Range("Q12").Value = 2.111978771320889
$ws.Range("R12").Value = 19.007808941888
$ws.Range("S12").Value = 0.01239475423611331
$ws.Range("T12").Value = 0.01239475423611331
$ws.Range("G13").Value = 2.406480666666666
$ws.Range("H13").Value = 7.219442
$ws.Range("I13").Value = 0.09799101667823351
$ws.Range("J13").Value = 0.09799101667823353
$ws.Range("M13").Value = 5.010059000000001
$ws.Range("N13").Value = 15.030177
$ws.Range("O13").Value = 0.7220833401192995
$ws.Range("P13").Value = 0.7220833401192994
$ws.Range("Q13").Value = 12.05661012235933
$ws.Range("R13").Value = 108.509491101234
$ws.Range("S13").Value = 0.07075768062470485
$ws.Range("T13").Value = 0.07075768062470485
$ws.Range("G14").Value = 2.383179
$ws.Range("H14").Value = 7.149537
$ws.Range("I14").Value = 0.09704218129443352
$ws.Range("J14").Value = 0.09704218129443352
$ws.Range("M14").Value = 1.050658666666667
$ws.Range("N14").Value = 3.151976
$ws.Range("O14").Value = 0.1514279810580986
$ws.Range("P14").Value = 0.1514279810580986
$ws.Range("Q14").Value = 2.503907670568001
$ws.Range("R14").Value = 22.53516903511201
$ws.Range("S14").Value = 0.01469490159089005
$ws.Range("T14").Value = 0.01469490159089005
$ws.Range("G15").Value = 2.383179
$ws.Range("H15").Value = 7.149537
$ws.Range("I15").Value = 0.09704218129443352
$ws.Range("J15").Value = 0.09704218129443352
$ws.Range("M15").Value = 0.8776213333333334
$ws.Range("N15").Value = 2.632864
$ws.Range("O15").Value = 0.126488678822602
$ws.Range("P15").Value = 0.126488678822602
$ws.Range("Q15").Value = 2.091528731552
$ws.Range("R15").Value = 18.823758583968
$ws.Range("S15").Value = 0.01227473730199631
$ws.Range("T15").Value = 0.01227473730199631
$ws.Range("G16").Value = 2.383179
$ws.Range("H16").Value = 7.149537
$ws.Range("I16").Value = 0.09704218129443352
$ws.Range("J16").Value = 0.09704218129443352
$ws.Range("M16").Value = 5.010059000000001
$ws.Range("N16").Value = 15.030177
$ws.Range("O16").Value = 0.7220833401192995
$ws.Range("P16").Value = 0.7220833401192994
$ws.Range("Q16").Value = 11.939867397561
$ws.Range("R16").Value = 107.458806578049
$ws.Range("S16").Value = 0.07007254240154717
$ws.Range("T16").Value = 0.07007254240154716
$ws.Range("G17").Value = 10.75644933333333
$ws.Range("H17").Value = 32.269348
$ws.Range("I17").Value = 0.4379987010164666
$ws.Range("J17").Value = 0.4379987010164666
$ws.Range("M17").Value = 1.050658666666667
$ws.Range("N17").Value = 3.151976
$ws.Range("O17").Value = 0.1514279810580986
$ws.Range("P17").Value = 0.1514279810580986
$ws.Range("Q17").Value = 11.30135671462756
$ws.Range("R17").Value = 101.712210431648
$ws.Range("S17").Value = 0.0663252590009933
$ws.Range("T17").Value = 0.06632525900099329
$ws.Range("G18").Value = 10.75644933333333
$ws.Range("H18").Value = 32.269348
$ws.Range("I18").Value = 0.4379987010164666
$ws.Range("J18").Value = 0.4379987010164666
$ws.Range("M18").Value = 0.8776213333333334
$ws.Range("N18").Value = 2.632864
$ws.Range("O18").Value = 0.126488678822602
$ws.Range("P18").Value = 0.126488678822602
$ws.Range("Q18").Value = 9.440089405852445
$ws.Range("R18").Value = 84.96080465267201
$ws.Range("S18").Value = 0.0554018770175887
$ws.Range("T18").Value = 0.0554018770175887
$ws.Range("G19").Value = 10.75644933333333
$ws.Range("H19").Value = 32.269348
$ws.Range("I19").Value = 0.4379987010164666
$ws.Range("J19").Value = 0.4379987010164666
$ws.Range("M19").Value = 5.010059000000001
$ws.Range("N19").Value = 15.030177
$ws.Range("O19").Value = 0.7220833401192995
$ws.Range("P19").Value = 0.7220833401192994
$ws.Range("Q19").Value = 53.89044579051068
$ws.Range("R19").Value = 485.0140121145961
$ws.Range("S19").Value = 0.3162715649978846
$ws.Range("T19").Value = 0.3162715649978846
